$d = $word.ActiveDocument

# Locate the "Module Usage" bullet that documents the CBX variant so we can
# find the following (currently empty) bullet that needs the new UBX text.
$anchor = $d.Content.Duplicate
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute("Document for more details", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'SetLO_CBX ... Document for more details' paragraph"
}

$cbxParagraph = $anchor.Paragraphs(1)
$targetParagraph = $cbxParagraph.Next()

# Sanity-check: the paragraph we are about to replace should be the empty
# numbered/bold bullet that immediately follows the CBX usage bullet.
if ($targetParagraph.Range.Text.Trim().Length -ne 0) {
    throw "Unexpected content in the paragraph that should be replaced"
}

$insertionPoint = $targetParagraph.Range.Duplicate
$insertionPoint.Collapse(0)

$newContentPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00357A04" w:rsidRPr="00357A04" w:rsidRDefault="00357A04" w:rsidP="00357A04"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">CreateMax2871Packets module(Sub module of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SetLO_UBX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) --&gt; Updated on 09-06-2018</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">. Check </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>SetLO_U</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>BX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Document for more details</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newContentPackage)

Write-Output "UBX usage bullet inserted."
